$d = $word.ActiveDocument

# Word's Find/Replace "smart quotes" autocorrect would turn straight apostrophes
# in a Replacement string into curly ones, so any snippet containing an
# apostrophe is applied via direct Range.Text assignment instead of
# Find.Execute's Replace parameter (which is fine for apostrophe-free snippets).

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Vertical-tab char maps to a <w:br/> line break inside a single run/paragraph
# when assigned through Range.Text (unlike "`r" which starts a new paragraph).
$lf = [char]11

# 1. Title heading: "Objet RS-RR:15-15:resourcesResponse" -> "resourcesResponse"
Replace-Text "Objet RS-RR:15-15:resourcesResponse" "resourcesResponse"

$t0 = $d.Tables(1)
$t1 = $d.Tables(2)

# 2. caseId description (Table 1, row 2, col 5): prepend "A valoriser avec "
#    and lowercase "Identifiant" -> "l'identifiant". Only the first <w:t> (up
#    to the first <w:br/>) changes; later runs must stay untouched, so locate
#    that exact first run's text and overwrite just that narrow range.
$rng = $d.Content
$rng.Find.Execute("Identifiant partagé de l'affaire/dossier, généré une seule fois par le système du partenaire qui recoit la primo-demande de secours (créateur du dossier). ") | Out-Null
$rng.Text = "A valoriser avec l'identifiant partagé de l'affaire/dossier, généré une seule fois par le système du partenaire qui recoit la primo-demande de secours (créateur du dossier). "

# 3. requestId description (Table 1, row 3, col 5): same pattern.
$rng2 = $d.Content
$rng2.Find.Execute("Identifiant unique partagé de la demande de ressource,  généré une seule fois par le système du partenaire qui émet la demande ") | Out-Null
$rng2.Text = "A valoriser avec l'identifiant unique partagé de la demande de ressource,  généré une seule fois par le système du partenaire qui émet la demande "

# 4. response row Description (Table 1, row 4, col 5)
$t0.Cell(4,5).Range.Text = "Objet permettant de transmettre les détails de la réponse à une demande de ressource"
# 5. response row Exemple (Table 1, row 4, col 6): clear it out
$t0.Cell(4,6).Range.Text = ""

# 6. Second heading: "Type response" -> "response"
Replace-Text "Type response" "response"

# 7. datetime row Description (Table 2, row 2, col 5)
$t1.Cell(2,5).Range.Text = "Groupe date heure de début de la demande"
# 8. datetime row Exemple (Table 2, row 2, col 6): was empty
$t1.Cell(2,6).Range.Text = "2022-09-27T08:23:34+02:00"

# 9. answer row Description (Table 2, row 3, col 5): two lines w/ <w:br/>
$t1.Cell(3,5).Range.Text = "A valoriser avec la réponse apportée. Cf Nomenclature associée" + $lf + "ACCEPTEE, REFUSEE, PARTIELLE, DIFFEREE"
# 10. answer row Exemple (Table 2, row 3, col 6): was empty
$t1.Cell(3,6).Range.Text = "ACCEPTEE"

# 11. deadline row Description (Table 2, row 4, col 5): apostrophes but single run
$t1.Cell(4,5).Range.Text = "A valoriser avec le délai de réponse auquel s'engage l'expéditeur (en minutes),"
# 12. deadline row Exemple (Table 2, row 4, col 6): was empty
$t1.Cell(4,6).Range.Text = "10"

# 13. freetext row Description (Table 2, row 5, col 5): apostrophe but single run
$t1.Cell(5,5).Range.Text = "Commentaire libre permettant d'apporter toutes précisions utiles à la réponse. Le motif de refus est notifié dans ce champ."
# 14. freetext row Exemple (Table 2, row 5, col 6): was empty
$t1.Cell(5,6).Range.Text = "SMUR 1 non dispo"
